# Update "想去人数" (want-to-go count) figures for a handful of events
# on the "展览" and "全部类型" sheets, reflecting refreshed numbers.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7819
$ws1.Range("F5").Value = 57
$ws1.Range("F6").Value = 573
$ws1.Range("F7").Value = 1188
$ws1.Range("F8").Value = 212
$ws1.Range("F10").Value = 174

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7819
$ws4.Range("F5").Value = 57
$ws4.Range("F6").Value = 573
$ws4.Range("F7").Value = 1188
$ws4.Range("F8").Value = 212
$ws4.Range("F11").Value = 174
